$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.970.79'
$ws.Range('E2').Value = '  -0.58%  '
$ws.Range('D3').Value = '1.635.37'
$ws.Range('E3').Value = '  -1.16%  '
$ws.Range('E4').Value = '  +0.15%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '212.25'
$ws.Range('E5').Value = '  -0.85%  '
$ws.Range('E6').Value = '  -1.02%  '
$ws.Range('E7').Value = '  +0.15%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '23.32'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.259'
$ws.Range('E9').Value = '  -2.82%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0615'
$ws.Range('E10').Value = '  -0.07%  '
$ws.Range('E11').Value = '  +1.06%  '
$ws.Range('E12').Value = '  -1.09%  '
$ws.Range('D13').Value = '1.635.08'
$ws.Range('E13').Value = '  -1.47%  '
$ws.Range('E14').Value = '  -0.64%  '
$ws.Range('E15').Value = '  -0.24%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.40'
$ws.Range('E16').Value = '  -0.61%  '
$ws.Range('D17').Value = '27.973.01'
$ws.Range('E17').Value = '  -0.45%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '231.06'
$ws.Range('E18').Value = '  -1.15%  '
$ws.Range('E19').Value = '  -0.16%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.55'
$ws.Range('E20').Value = '  -2.26%  '
$ws.Range('E21').Value = '  +0.17%  '
$ws.Range('E22').Value = '  -0.89%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.38'
$ws.Range('E23').Value = '  -2.91%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.07'
$ws.Range('E24').Value = '  -4.14%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '155.03'
$ws.Range('E25').Value = '  +1.77%  '
$ws.Range('E26').Value = '  +0.55%  '
$ws.Range('B27').Value = 'Stellar'
$ws.Range('C27').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.111'
$ws.Range('E27').Value = '  -0.59%  '
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.67'
$ws.Range('E28').Value = '  -0.96%  '
$ws.Range('E29').Value = '  +0.15%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.19'
$ws.Range('E30').Value = '  -0.54%  '
$ws.Range('E31').Value = '  -0.28%  '
$ws.Range('E32').Value = '  +1.05%  '
$ws.Range('D33').Value = '1.408.18'
$ws.Range('E33').Value = '  -3.14%  '
$ws.Range('E34').Value = '  -0.37%  '
$ws.Range('E35').Value = '  -0.09%  '
$ws.Range('E36').Value = '  +9.50%  '
$ws.Range('E37').Value = '  +1.38%  '
$ws.Range('E38').Value = '  +0.28%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.563'
$ws.Range('E39').Value = '  +0.25%  '
$ws.Range('E40').Value = '  -2.27%  '
$ws.Range('E41').Value = '  +0.10%  '
$ws.Range('E42').Value = '  +0.10%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '67.02'
$ws.Range('E43').Value = '  -3.66%  '
$ws.Range('E44').Value = '  +2.32%  '
$ws.Range('E45').Value = '  +0.04%  '
$ws.Range('E46').Value = '  -4.12%  '
$ws.Range('D47').Value = '1.776.61'
$ws.Range('E47').Value = '  -1.07%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '87.89'
$ws.Range('E48').Value = '  -1.32%  '
$ws.Range('D49').Value = '0.0₆0106'
$ws.Range('E49').Value = '  +8.74%  '
$ws.Range('E50').Value = '  -1.30%  '
$ws.Range('E51').Value = '  -0.33%  '
